# Apply "Added legends to corpuses" edit:
#  1. Rename Sheet1 -> Corpus
#  2. Add a new sheet "Legend" with a two-column key/description table
#  3. Make "Legend" the active sheet/tab, clear the Corpus sheet's
#     top-left/selection scroll state

$wb = $excel.ActiveWorkbook

# --- 1. Rename the first sheet -------------------------------------------------
$wsCorpus = $wb.Worksheets.Item(1)
$wsCorpus.Name = "Corpus"

# --- 2. Add the Legend sheet ----------------------------------------------------
# Add(Before, After) -> placing it After the Corpus sheet makes it the 2nd tab
$wsLegend = $wb.Worksheets.Add($null, $wsCorpus)
$wsLegend.Name = "Legend"

$legend = @(
    @("id", "Id of the entry in the corpus"),
    @("name", "mention label as it appears in the original text"),
    @("doi", "doi of the paper the mention name appears"),
    @("papragraph", "text where the original mention appears in the paper"),
    @("field/topic/keywords", "topics of the paper as obtained from OpenALEX. Only those with confidence >0.5 are considered"),
    @("authors", "authors of the paper"),
    @("authors_oa", "authors of the paper (Open Alex ids)"),
    @("url (groung truth)", "URL that the annotators think the software refers to. "),
    @("annotator", "Person or persons responsible for validating the ground truth"),
    @("comments", "Possible comments about the mention"),
    @("candidate_urls", "List of possible URLs that software refers to"),
    @("language", "Language fetched from paragraph is exists ")
)

for ($i = 0; $i -lt $legend.Count; $i++) {
    $row = $i + 1
    $wsLegend.Cells.Item($row, 1).Value = $legend[$i][0]
    $wsLegend.Cells.Item($row, 2).Value = $legend[$i][1]
}

# --- 3. Formatting: shrink the legend text to 10pt ------------------------------
# Most of the sheet (A1:B10 + A12) got only a size tweak; A11 was left at the
# sheet default; B11/B12 additionally had the font name re-applied explicitly.
$wsLegend.Range("A1:B10").Font.Size = 10
$wsLegend.Range("A12").Font.Size = 10

$wsLegend.Range("B11:B12").Font.Size = 10
$wsLegend.Range("B11:B12").Font.Name = "Calibri"

# --- 4. Selections / view state -------------------------------------------------
$wsCorpus.Range("I471").Select()
$wsLegend.Range("A13").Select()
$wsLegend.Activate()
